$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Make the settings sheet the active / tab-selected sheet (moves tabSelected
# from "survey" to "settings" and sets activeTab on the workbook view).
$ws.Activate()

# Add the new "disableSwipeNavigation" setting row.
$ws.Range("A6").Value = "disableSwipeNavigation"
$ws.Range("B6").Value = $true

# Formatting for the new cells: Times New Roman 12pt black for the name,
# Arial 10pt black for the boolean value, both wrapped.
$ws.Range("A6").Font.Name = "Times New Roman"
$ws.Range("A6").Font.Color = 0
$ws.Range("A6").WrapText = $true

$ws.Range("B6").Font.Name = "Arial"
$ws.Range("B6").Font.Size = 10
$ws.Range("B6").Font.Color = 0
$ws.Range("B6").WrapText = $true

# Select the new row, with A6 as the active cell.
$ws.Range("A6:B6").Select()
